# Auto-generated edit script: update cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "87.031.02"
$ws.Range("E2").Value = "  -3.69%  "

# Row 3
$ws.Range("D3").Value = "3.051.02"
$ws.Range("E3").Value = "  -4.43%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.358"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -11.17%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.767"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.90%  "

# Row 9
$ws.Range("E9").Value = "  +0.20%  "

# Row 10
$ws.Range("D10").Value = "3.048.46"
$ws.Range("E10").Value = "  -4.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.575"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.91%  "

# Row 12
$ws.Range("E12").Value = "  -0.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -10.41%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.38%  "

# Row 15
$ws.Range("D15").Value = "87.042.98"
$ws.Range("E15").Value = "  -3.40%  "

# Row 16
$ws.Range("D16").Value = "3.624.92"
$ws.Range("E16").Value = "  -4.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.94%  "

# Row 18
$ws.Range("D18").Value = "3.065.03"
$ws.Range("E18").Value = "  -3.62%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000206"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "412.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.40%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.80%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.46%  "

# Row 28
$ws.Range("D28").Value = "3.228.58"
$ws.Range("E28").Value = "  -3.59%  "

# Row 29
$ws.Range("E29").Value = "  +0.12%  "

# Row 30
$ws.Range("E30").Value = "  +0.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.148"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.73%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "489.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -14.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.139"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.15%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.36%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.05%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.74%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.04%  "

# Row 41
$ws.Range("E41").Value = "  +0.33%  "

# Row 42
$ws.Range("E42").Value = "  -0.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.355"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.32%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "146.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.44%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.129"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.98%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.77%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0635"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.81%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "156.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.697"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.33%  "
